$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing values in column D (rows 129-132)
$ws.Range("D129").Value = 0.7136611031493167
$ws.Range("D130").Value = 0.7250871391493168
$ws.Range("D131").Value = 0.6447426901493167
$ws.Range("D132").Value = 0.7064651831493167

# Update existing values in column C (rows 133-139)
$ws.Range("C133").Value = 0.5279710641493167
$ws.Range("C134").Value = -0.2139369238506833
$ws.Range("C135").Value = 0.1647551661493167
$ws.Range("C136").Value = 0.4518305101493167
$ws.Range("C137").Value = 0.4794748351493167
$ws.Range("C138").Value = 0.2386249091493167
$ws.Range("C139").Value = 0.3244906151493167

# Add new values in column D (rows 136-139)
$ws.Range("D136").Value = 0.791995474
$ws.Range("D137").Value = 0.788120887
$ws.Range("D138").Value = 0.597740902
$ws.Range("D139").Value = 0.620527487

# Update existing values in column B (rows 140-141)
$ws.Range("B140").Value = 0.03316543414931669
$ws.Range("B141").Value = -0.0107480648506833

# Add new values in column C (rows 140-145)
$ws.Range("C140").Value = 0.241887844
$ws.Range("C141").Value = 0.331651578
$ws.Range("C142").Value = 0.154182215
$ws.Range("C143").Value = 0.166899468
$ws.Range("C144").Value = 0.042359665
$ws.Range("C145").Value = 0.266698307
